$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F6").Value = 0.798095238095238
$ws.Range("G6").Value = 0.7678412698412698
$ws.Range("H6").Value = 0.798095238095238
$ws.Range("I6").Value = 0.756998556998557

$ws.Range("N6").Value = 0.7704761904761905
$ws.Range("O6").Value = 0.7492316017316017
$ws.Range("P6").Value = 0.7704761904761905
$ws.Range("Q6").Value = 0.7255743733390793
$ws.Range("R6").Value = 0.7304761904761905
$ws.Range("S6").Value = 0.7195346320346321
$ws.Range("T6").Value = 0.7304761904761905
$ws.Range("U6").Value = 0.6792513107807225
$ws.Range("V6").Value = 0.7295238095238095
$ws.Range("W6").Value = 0.7086349206349206
$ws.Range("X6").Value = 0.7295238095238095
$ws.Range("Y6").Value = 0.6783219954648526
